$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 18 values ---
$ws.Range("C18").Value = 26394
$ws.Range("K18").Value = 58430
$ws.Range("L18").Value = 26396.6394
$ws.Range("T18").Value = 58435.843

# --- Add new rows 19-23 ---
$ws.Range("A19").Value = "'2025-01-08"
$ws.Range("B19").Value = 18
$ws.Range("C19").Value = 24734
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 16117
$ws.Range("F19").Value = 15506
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 56357
$ws.Range("L19").Value = 24736.4734
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 16118.6117
$ws.Range("O19").Value = 15507.5506
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("S19").Value = 0
$ws.Range("T19").Value = 56362.6357
$ws.Range("A19").Style = "Normal"

$ws.Range("A20").Value = "'2025-01-08"
$ws.Range("B20").Value = 19
$ws.Range("C20").Value = 24108
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 15648
$ws.Range("F20").Value = 15435
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 55191
$ws.Range("L20").Value = 24110.4108
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 15649.5648
$ws.Range("O20").Value = 15436.5435
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("S20").Value = 0
$ws.Range("T20").Value = 55196.5191
$ws.Range("A20").Style = "Normal"

$ws.Range("A21").Value = "'2025-01-08"
$ws.Range("B21").Value = 20
$ws.Range("C21").Value = 23135
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 14515
$ws.Range("F21").Value = 14707
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 52357
$ws.Range("L21").Value = 23137.3135
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 14516.4515
$ws.Range("O21").Value = 14708.4707
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0
$ws.Range("R21").Value = 0
$ws.Range("S21").Value = 0
$ws.Range("T21").Value = 52362.2357
$ws.Range("A21").Style = "Normal"

$ws.Range("A22").Value = "'2025-01-08"
$ws.Range("B22").Value = 21
$ws.Range("C22").Value = 22359
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 13506
$ws.Range("F22").Value = 14749
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 50614
$ws.Range("L22").Value = 22361.2359
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 13507.3506
$ws.Range("O22").Value = 14750.4749
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0
$ws.Range("R22").Value = 0
$ws.Range("S22").Value = 0
$ws.Range("T22").Value = 50619.0614
$ws.Range("A22").Style = "Normal"

$ws.Range("A23").Value = "'2025-01-08"
$ws.Range("B23").Value = 22
$ws.Range("C23").Value = 22099
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 11221
$ws.Range("F23").Value = 14631
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 47951
$ws.Range("L23").Value = 22101.2099
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 11222.1221
$ws.Range("O23").Value = 14632.4631
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0
$ws.Range("R23").Value = 0
$ws.Range("S23").Value = 0
$ws.Range("T23").Value = 47955.7951
$ws.Range("A23").Style = "Normal"

